$d = $word.ActiveDocument

# --- Bug fix: merge the "Requirements" / ":" runs into a single run ---
# A find/replace over the whole (already-correct) text collapses the two
# adjacent runs that spell out "Requirements" + ":" into one run holding
# "Requirements:" (matching the canonical OOXML).
$d.Content.Find.Execute("Requirements:", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Requirements:", 2) | Out-Null

# --- Add the missing "Multi-User Login" child item under "Update Task" ---
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "Update Task") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $newRange = $target.Range.InsertParagraphAfter()
    $target.Next().Range.Text = "Multi-User Login"
}
